$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column C (rows 2-402) hot water intensity values to 0.05
$ws.Range("C2:C402").Value = 0.05

# Update the active cell selection to C12
$ws.Range("C12").Select()
